$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "总计" (summary) sheet: a new 2022-Q4 row is inserted at the top
#    of the data (row 2), pushing every existing row down by one.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# Copy the number-style (bold/centered, same as the other index cells)
# from A3 onto the freshly inserted A2, then strip the stray formatting
# that Insert() propagated into B2:D2 so they stay "plain" like the rest
# of the data rows.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()
$excel.CutCopyMode = 0

# New first data row: 2022-Q4
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.07000000000000001

# Re-number the index column + keep the rest of the (shifted) rows in sync
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 5
$summary.Range("D4").Value = 0.12

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.16

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 0.53

# ------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet right after "总计" holding the
#    per-fund breakdown for the new quarter (same layout/style as the
#    other quarter sheets).
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Header row (bold, centered - matches the other quarter sheets)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$q4.Range("B1:H1").Font.Bold = $true
$q4.Range("B1:H1").HorizontalAlignment = -4108
$q4.Range("B1:H1").VerticalAlignment = -4160

# Index column (A) is numeric/bold/centered like the other sheets
$q4.Range("A2:A3").Font.Bold = $true
$q4.Range("A2:A3").HorizontalAlignment = -4108
$q4.Range("A2:A3").VerticalAlignment = -4160
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1

# Fund-code / numeric-looking metric columns are stored as TEXT, same as
# every other quarter sheet in this workbook.
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "160518"
$q4.Range("C2").Value = "博时睿远事件驱动灵活配置混合（LOF）"
$q4.Range("D2").Value = "1.79"
$q4.Range("E2").Value = "83.89"
$q4.Range("F2").Value = "2.62"
$q4.Range("G2").Value = "0.0469"
$q4.Range("H2").Value = 9

$q4.Range("B3").Value = "159804"
$q4.Range("C3").Value = "国寿安保国证创业板中盘精选88ETF"
$q4.Range("D3").Value = "1.15"
$q4.Range("E3").Value = "99.00"
$q4.Range("F3").Value = "1.83"
$q4.Range("G3").Value = "0.0210"
$q4.Range("H3").Value = 8
